$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook's styles.xml ships with no cellXfs entries, which the engine
# needs initialized before any cell write; touching a named Style first
# populates it without altering any existing cell content.
$ws.Cells.Style = "Normal"

# Insert a new row above the current row 1, pushing the existing data
# (1/bee/3.14 and 2/butterfly/6.28) down to rows 2 and 3.
$ws.Rows.Item(1).Insert()

# Fill in the new header row.
$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"
$ws.Range("C1").Value = "c"
